$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.652167666666667
$ws.Range("H2").Value = 16.956503
$ws.Range("I2").Value = 0.1860329065948871
$ws.Range("J2").Value = 0.1860329065948871
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 119.1579455614248
$ws.Range("R2").Value = 1072.421510052823
$ws.Range("S2").Value = 0.01063375283228868
$ws.Range("T2").Value = 0.01063375283228868

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.652167666666667
$ws.Range("H3").Value = 16.956503
$ws.Range("I3").Value = 0.1860329065948871
$ws.Range("J3").Value = 0.1860329065948871
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 1704.694608555593
$ws.Range("R3").Value = 15342.25147700034
$ws.Range("S3").Value = 0.1521283455879224
$ws.Range("T3").Value = 0.1521283455879224

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.652167666666667
$ws.Range("H4").Value = 16.956503
$ws.Range("I4").Value = 0.1860329065948871
$ws.Range("J4").Value = 0.1860329065948871
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 260.7641664595275
$ws.Range("R4").Value = 2346.877498135747
$ws.Range("S4").Value = 0.02327080817467596
$ws.Range("T4").Value = 0.02327080817467596

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.407289
$ws.Range("H5").Value = 49.221867
$ws.Range("I5").Value = 0.5400221369958743
$ws.Range("J5").Value = 0.5400221369958743
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 345.8954094731497
$ws.Range("R5").Value = 3113.058685258347
$ws.Range("S5").Value = 0.03086799015232012
$ws.Range("T5").Value = 0.03086799015232012

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.407289
$ws.Range("H6").Value = 49.221867
$ws.Range("I6").Value = 0.5400221369958743
$ws.Range("J6").Value = 0.5400221369958743
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 4948.440801616965
$ws.Range("R6").Value = 44535.96721455269
$ws.Range("S6").Value = 0.4416029173856635
$ws.Range("T6").Value = 0.4416029173856635

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.407289
$ws.Range("H7").Value = 49.221867
$ws.Range("I7").Value = 0.5400221369958743
$ws.Range("J7").Value = 0.5400221369958743
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 756.9543743681539
$ws.Range("R7").Value = 6812.589369313384
$ws.Range("S7").Value = 0.06755122945789074
$ws.Range("T7").Value = 0.06755122945789074

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.323166333333333
$ws.Range("H8").Value = 24.969499
$ws.Range("I8").Value = 0.2739449564092387
$ws.Range("J8").Value = 0.2739449564092387
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 175.4674417560065
$ws.Range("R8").Value = 1579.206975804059
$ws.Range("S8").Value = 0.01565885847524571
$ws.Range("T8").Value = 0.01565885847524571

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.323166333333333
$ws.Range("H9").Value = 24.969499
$ws.Range("I9").Value = 0.2739449564092387
$ws.Range("J9").Value = 0.2739449564092387
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 2510.268203510728
$ws.Range("R9").Value = 22592.41383159655
$ws.Range("S9").Value = 0.2240183941835933
$ws.Range("T9").Value = 0.2240183941835933

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.323166333333333
$ws.Range("H10").Value = 24.969499
$ws.Range("I10").Value = 0.2739449564092387
$ws.Range("J10").Value = 0.2739449564092387
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 383.9913568055279
$ws.Range("R10").Value = 3455.922211249751
$ws.Range("S10").Value = 0.03426770375039966
$ws.Range("T10").Value = 0.03426770375039966

